$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 = Shopkeeper (商人) card: append the new repeatable discard-for-gold option (③)
$ws.Range("D9").Value = "多选：①可重复 支付1金币为1张道具牌充1能。②支付3金币，从遗物牌堆翻开3张牌，选其中1张获得。③可重复 弃置1张任意战利品牌，获得1金币。"

# Update the view to match where the author left the selection/scroll position
$ws.Range("C9").Select()
$excel.ActiveWindow.ScrollRow = 3
